$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently occupies rows 4-16 (with rows 4 being header).
# The edit shifts the whole table up so it starts at row 1 (rows 1-13),
# and fills previously-blank percentage cells with explicit 0 values.

# Delete rows 1-3 so the existing table (rows 4-16) shifts up to rows 1-13.
$ws.Range("A1:A3").EntireRow.Delete() | Out-Null

# Fill in the previously blank cells (Support, Business Development, Training rows)
# with explicit 0 values. After the shift these are rows 5, 7, and 9.
$zeroRows = @(5, 7, 9)
foreach ($r in $zeroRows) {
    for ($c = 2; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# Update the sheet view: move the active selection to I12 (this also clears
# the previous scrolled-down view position).
$ws.Range("I12").Select()
